$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 348, shifting rows 348:415 down to 349:416
$ws.Rows("348").Insert()

# Populate the newly inserted row 348 with the new record
$ws.Range("A348").Value = 5
$ws.Range("B348").Value = "Macroferia Regional de Talca"
$ws.Range("C348").Value = "Maule"
$ws.Range("D348").Value = 44694
$ws.Range("E348").Value = 7
$ws.Range("F348").Value = 100112043
$ws.Range("G348").Value = "Pepino ensalada"
$ws.Range("H348").Value = "Sin especificar"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 300
$ws.Range("K348").Value = 15000
$ws.Range("L348").Value = 15000
$ws.Range("M348").Value = 15000
$ws.Range("N348").Value = '$/caja 60 unidades'
$ws.Range("O348").Value = "Región de Arica y Parinacota"
$ws.Range("P348").Value = 250
$ws.Range("Q348").Value = 60
$ws.Range("R348").Value = "Hortaliza"
